$d = $word.ActiveDocument

# The document currently holds a single empty paragraph. Replace its
# range with the fully-authored methods paragraph (text + run/paragraph
# formatting + spell-check proof markers) using the paragraph's own Range.
$p = $d.Paragraphs(1)
$r = $p.Range

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:before="180" w:after="180"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t>Focal follow data was collected at two locations at Palmyra Atoll National Wildlife Refuge, Penguin Spit and Western Terrace. Individual fish were identified by unique markings on the face, tail, and body as well as missing scales or scars.  A diver would identify an individual fish and follow the fish for 2-3mins before stating the focal follow. The diver would tow a GPS that was recording a position every 5 secs and synchronized the time on their wrist watch to the GPS time.  The diver would then follow the fish and record the start and stop time of each activity, while also following the path of the fish.  GPS tracks were downloaded and each position was associated with an activity and activity summaries were calculated.  For</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
      <w:i/>
      <w:iCs/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t xml:space="preserve"> C. </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
      <w:i/>
      <w:iCs/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t>Microrhinos</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t xml:space="preserve"> we calculated the 95% kernel utilization distributions (KUD) for the entire focal follow as well as for only the locations categorized as ‘feeding’.  All KUD estimates were done in R with the </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t>adehabitatHR</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t xml:space="preserve"> package. </w:t>
  </w:r>
</w:p>
"@

[void]$r.InsertXML($xml)
